# Junction_Flooding_478.xlsx update
# - custom accuracy column widths (several 7 -> 8, one 8 -> 9)
# - replace data rows 2-5 with new simulation values ("데이터 1000개")
# - remove old row 6 (sheet now spans A1:AH5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update column widths per diff -------------------------------------
# ColumnWidth (Excel character units) = target XML width - 0.8333333333333333
$colWidths = @{
    2  = 8
    3  = 8
    7  = 8
    9  = 8
    10 = 8
    11 = 8
    12 = 8
    15 = 8
    16 = 8
    17 = 8
    20 = 9
    22 = 8
    23 = 8
    24 = 8
    26 = 8
    27 = 8
    28 = 8
    29 = 8
    30 = 8
    34 = 8
}
$offset = 0.8333333333333333
foreach ($colIndex in $colWidths.Keys) {
    $ws.Columns($colIndex).ColumnWidth = $colWidths[$colIndex] - $offset
}

# --- 2. Overwrite data rows 2-5 with the new readings ----------------------
$row2 = @(45159.50694444445,14.835,9.791,3.698,32.243,24.166,11.51,34.958,18.033,7.29,10.735,12.533,13.25,3.739,11.655,16.06,10.282,3.096,1.74,170.025,32.298,10.758,20.812,10.713,2.837,18.288,9.502,8.642,10.303,12.679,3.311,31.418,5.68,13.449)
$row3 = @(45159.51388888889,24.448,17.766,2.046,53.468,42.836,19.122,71.922,29.668,12.979,19.105,21.28,22.55,6.159,19.174,27.132,16.371,1.532,1.262,284.539,53.597,17.698,35.728,18.695,3.139,35.631,15.633,13.955,16.449,22.082,1.246,65.728,9.824,22.127)
$row4 = @(45159.52083333334,1.392,0.676,0.806,3.287,1.764,1.012,13.319,1.745,0.744,0.601,1.215,1.336,0.402,1.128,1.664,1.299,0.878,0.369,10.02,3.807,1.041,2.412,1.065,0.567,5.777,0.92,1.003,1.162,1.069,0.766,13.035,0.402,1.316)
$row5 = @(45159.52777777778,18.21,13.42,1.16,39.84,32.21,14.26,51.32,22.11,9.75,14.39,15.91,16.89,4.59,14.29,20.25,12.12,0.8,0.77,210.11,39.79,13.19,26.65,14.02,2.19,25.38,11.65,10.36,12.19,16.57,0.55,46.33,7.37,16.49)

$dataRows = @{
    2 = $row2
    3 = $row3
    4 = $row4
    5 = $row5
}

foreach ($r in $dataRows.Keys) {
    $vals = $dataRows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

# --- 3. Drop the old row 6 (data set now has one fewer row) ---------------
$ws.Rows(6).Delete()
